$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4224.75
$ws.Range("I64").Value = 3966.3333
$ws.Range("K64").Value = 3966.3333
$ws.Range("M64").Value = -3718.3333
$ws.Range("H67").Value = 4224.75
$ws.Range("I67").Value = 3966.3333
$ws.Range("K67").Value = 3966.3333
$ws.Range("M67").Value = -3108.3333
$ws.Range("H69").Value = 6500
$ws.Range("J69").Value = 6500
$ws.Range("L69").Value = 19500
$ws.Range("N69").Value = -21248
$ws.Range("H72").Value = 6500
$ws.Range("J72").Value = 6500
$ws.Range("L72").Value = 58500
$ws.Range("N72").Value = -67236
$ws.Range("H127").Value = 2511.25
$ws.Range("I127").Value = 1285
$ws.Range("J127").Value = 2817.8125
$ws.Range("K127").Value = 3855
$ws.Range("L127").Value = 8453.4375
$ws.Range("M127").Value = 1105
$ws.Range("N127").Value = -18373.4375
$ws.Range("H129").Value = 824008
$ws.Range("I129").Value = 372.3
$ws.Range("J129").Value = 1059332.5
$ws.Range("K129").Value = 1116.9
$ws.Range("L129").Value = 3177997.5
$ws.Range("M129").Value = 3883.1
$ws.Range("N129").Value = -3187997.5
$ws.Range("H137").Value = 789.57446
$ws.Range("I137").Value = 717.8946999999999
$ws.Range("J137").Value = 1092.2222
$ws.Range("K137").Value = 2153.6841
$ws.Range("L137").Value = 3276.6666
$ws.Range("M137").Value = 396.3159000000001
$ws.Range("N137").Value = -8376.6666
$ws.Range("H138").Value = 2431.1685
$ws.Range("I138").Value = 1081.3256
$ws.Range("J138").Value = 3547.3845
$ws.Range("K138").Value = 3243.976799999999
$ws.Range("L138").Value = 10642.1535
$ws.Range("M138").Value = 1896.023200000001
$ws.Range("N138").Value = -20922.1535

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1398.1936
$ws.Range("I45").Value = 1388.3636
$ws.Range("J45").Value = 1422.2222
$ws.Range("K45").Value = 1388.3636
$ws.Range("L45").Value = 1422.2222
$ws.Range("M45").Value = -1011.3636
$ws.Range("N45").Value = -2176.2222
$ws.Range("H74").Value = 1260.125
$ws.Range("I74").Value = 1403.6842
$ws.Range("K74").Value = 1403.6842
$ws.Range("M74").Value = -529.6841999999999
$ws.Range("H77").Value = 1260.125
$ws.Range("I77").Value = 1403.6842
$ws.Range("K77").Value = 7018.420999999999
$ws.Range("M77").Value = -2650.420999999999
$ws.Range("H140").Value = 59227
$ws.Range("J140").Value = 59227
$ws.Range("L140").Value = 59227
$ws.Range("N140").Value = -69587

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 6020
$ws.Range("J74").Value = 6020
$ws.Range("L74").Value = 6020
$ws.Range("N74").Value = -7892
$ws.Range("H77").Value = 6020
$ws.Range("J77").Value = 6020
$ws.Range("L77").Value = 18060
$ws.Range("N77").Value = -27420

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5369.8213
$ws.Range("I58").Value = 2069.5
$ws.Range("J58").Value = 9770.25
$ws.Range("K58").Value = 2069.5
$ws.Range("L58").Value = 9770.25
$ws.Range("M58").Value = -1866.5
$ws.Range("N58").Value = -10176.25
$ws.Range("H132").Value = 1511.6
$ws.Range("I132").Value = 1124.1111
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3372.3333
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -842.3333000000002
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 5369.8213
$ws.Range("I136").Value = 2069.5
$ws.Range("J136").Value = 9770.25
$ws.Range("K136").Value = 6208.5
$ws.Range("L136").Value = 29310.75
$ws.Range("M136").Value = -3658.5
$ws.Range("N136").Value = -34410.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 867.0833
$ws.Range("J113").Value = 648
$ws.Range("L113").Value = 1944
$ws.Range("N113").Value = -6284
$ws.Range("H115").Value = 1277.8182
$ws.Range("I115").Value = 495.2
$ws.Range("J115").Value = 1930
$ws.Range("K115").Value = 1485.6
$ws.Range("L115").Value = 5790
$ws.Range("M115").Value = -310.5999999999999
$ws.Range("N115").Value = -8140
$ws.Range("H118").Value = 3709.8462
$ws.Range("I118").Value = 421.33334
$ws.Range("J118").Value = 6528.5713
$ws.Range("K118").Value = 1264.00002
$ws.Range("L118").Value = 19585.7139
$ws.Range("M118").Value = -21.00001999999995
$ws.Range("N118").Value = -22071.7139
$ws.Range("H122").Value = 1060.386
$ws.Range("I122").Value = 560.625
$ws.Range("J122").Value = 1141.9796
$ws.Range("K122").Value = 5045.625
$ws.Range("L122").Value = 10277.8164
$ws.Range("M122").Value = -2595.625
$ws.Range("N122").Value = -15177.8164

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3985.0322
$ws.Range("I70").Value = 3922.5908
$ws.Range("J70").Value = 4137.6665
$ws.Range("K70").Value = 3922.5908
$ws.Range("L70").Value = 4137.6665
$ws.Range("M70").Value = -3652.5908
$ws.Range("N70").Value = -4677.6665
$ws.Range("H73").Value = 3985.0322
$ws.Range("I73").Value = 3922.5908
$ws.Range("J73").Value = 4137.6665
$ws.Range("K73").Value = 3922.5908
$ws.Range("L73").Value = 4137.6665
$ws.Range("M73").Value = -2986.5908
$ws.Range("N73").Value = -6009.6665
$ws.Range("H80").Value = 3498.111
$ws.Range("I80").Value = 3471.6956
$ws.Range("K80").Value = 3471.6956
$ws.Range("M80").Value = -2473.6956
$ws.Range("H83").Value = 3498.111
$ws.Range("I83").Value = 3471.6956
$ws.Range("K83").Value = 17358.478
$ws.Range("M83").Value = -12366.478
$ws.Range("H132").Value = 5154.1113
$ws.Range("I132").Value = 6193.778
$ws.Range("J132").Value = 4114.4443
$ws.Range("K132").Value = 18581.334
$ws.Range("L132").Value = 12343.3329
$ws.Range("M132").Value = -16051.334
$ws.Range("N132").Value = -17403.3329
$ws.Range("H134").Value = 12995.143
$ws.Range("J134").Value = 12995.143
$ws.Range("L134").Value = 38985.429
$ws.Range("N134").Value = -44055.429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2852.3157
$ws.Range("I68").Value = 2774.5
$ws.Range("J68").Value = 2985.7144
$ws.Range("K68").Value = 2774.5
$ws.Range("L68").Value = 2985.7144
$ws.Range("M68").Value = -2025.5
$ws.Range("N68").Value = -4483.7144
$ws.Range("H71").Value = 2852.3157
$ws.Range("I71").Value = 2774.5
$ws.Range("J71").Value = 2985.7144
$ws.Range("K71").Value = 13872.5
$ws.Range("L71").Value = 14928.572
$ws.Range("M71").Value = -10128.5
$ws.Range("N71").Value = -22416.572
$ws.Range("H135").Value = 34189.5
$ws.Range("J135").Value = 34189.5
$ws.Range("L135").Value = 34189.5
$ws.Range("N135").Value = -44329.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6586.6665
$ws.Range("I62").Value = 7250
$ws.Range("J62").Value = 6484.615
$ws.Range("K62").Value = 7250
$ws.Range("L62").Value = 6484.615
$ws.Range("M62").Value = -6626
$ws.Range("N62").Value = -7732.615
$ws.Range("H65").Value = 6586.6665
$ws.Range("I65").Value = 7250
$ws.Range("J65").Value = 6484.615
$ws.Range("K65").Value = 36250
$ws.Range("L65").Value = 32423.075
$ws.Range("M65").Value = -33130
$ws.Range("N65").Value = -38663.075
$ws.Range("H122").Value = 50000990
$ws.Range("I122").Value = 62500970
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 187502910
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -187500460
$ws.Range("N122").Value = -8200
